# Update the "database search with infinite scroll" data: refresh the
# timestamp column with a new run (2015-11-18) and extend the log from
# 5 data rows (3-7) to 12 data rows (3-14), updating the temperature
# reading from 20.79 to 18.44 for every pump column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New descending timestamp sequence for column A, rows 3..14.
$timestamps = @(
    "2015-11-18 11:40",
    "2015-11-18 11:30",
    "2015-11-18 11:20",
    "2015-11-18 11:10",
    "2015-11-18 11:00",
    "2015-11-18 10:50",
    "2015-11-18 10:40",
    "2015-11-18 10:30",
    "2015-11-18 10:20",
    "2015-11-18 10:10",
    "2015-11-18 10:00",
    "2015-11-18 09:40"
)

$newTemp = 18.44

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = 3 + $i

    # Timestamp column.
    $ws.Cells.Item($row, 1).Value = $timestamps[$i]

    # Pump blocks: (B,C,D) / (E,F,G) / (H,I,J) / (K,L,M) each = status, temp, mode.
    foreach ($col0 in @(2, 5, 8, 11)) {
        $ws.Cells.Item($row, $col0).Value = "OFF"
        $ws.Cells.Item($row, $col0 + 1).Value = $newTemp
        $ws.Cells.Item($row, $col0 + 2).Value = "자동"
    }
}

# Copy the formatting of the last original data row down onto the newly
# added rows so the new rows match the existing table look (style index 3
# on every cell, same as rows 3-7).
$ws.Range("A7:M7").Copy()
$ws.Range("A8:M14").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
